$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All features")

# --- New "library" feature block (rows 33-35), mirrors the B3:C5 pattern ---
# Row 33: category header ("категория" / "NEW (value-driven)") + right-hand table row 31
$ws.Range("B33").Value = "категория"
$ws.Range("C33").Value = "NEW (value-driven)"

# Row 34: "признаки" / full feature description
$ws.Range("B34").Value = "признаки"
$ws.Range("C34").Value = "Число библиотек - library (шт.) (8017001)"

# Row 35: blank continuation row (matches the blank B5/C5-style trailing row)
$ws.Range("B35").Value = ""
$ws.Range("C35").Value = ""

# Right-hand "Название / Кол-во записей / Диапазон" table, row 31 (sheet row 33)
$ws.Range("G33").Value = 31
$ws.Range("H33").Value = "library"
$ws.Range("I33").Value = 195097
$ws.Range("J33").Value = "2006 - 2017"

# Copy the row-above look (fill + centered alignment) onto the new H33 cell,
# then make sure the cell keeps the plain (non-bold/non-special) font.
$ws.Range("H30").Copy()
$ws.Range("H33").PasteSpecial(-4122)

# Reuse styling already present on the sheet for the new rows so the look
# matches the rest of the feature list / right-hand table.
$ws.Range("B3").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C33").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C34").PasteSpecial(-4122)

$ws.Range("B5").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C35").PasteSpecial(-4122)

$ws.Range("G32").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("I32").Copy()
$ws.Range("I33").PasteSpecial(-4122)
$ws.Range("J31").Copy()
$ws.Range("J33").PasteSpecial(-4122)

# --- Activate "All features" as the selected/visible sheet (it becomes the
# one the workbook was saved with focused), and move the selection to J37.
$ws.Activate()
$ws.Range("J37").Select()

$sortWs = $wb.Worksheets.Item("sort")
$sortWs.Range("I18").Select()
$ws.Activate()
